$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C on these rows holds numeric-looking attendance/minutes figures that
# are stored as text (shared strings) in the workbook, not real numbers.
# Force a text number format before writing so Excel doesn't silently
# reinterpret "66" etc. as a number, then clear the format again so we don't
# leave a stray style behind.
$ws.Range("C1:C7").NumberFormat = "@"

$ws.Range("C1").Value = "66"

$ws.Range("B2").Value = "February 26th 2022"
$ws.Range("C2").Value = "52"

$ws.Range("C3").Value = "71"

$ws.Range("C4").Value = "50"

$ws.Range("C5").Value = "83"

$ws.Range("C6").Value = "76"

$ws.Range("C7").Value = "84"

$ws.Range("C1:C7").ClearFormats()
